# Refresh market-price-derived Leve profit columns (H-N) with latest Universalis data.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* cells for the rows whose market
# data changed. A couple of rows lose all price data (H/J/L collapse to 0) and, to mirror
# the source sheet, their LeveProfitHQ (N) cell is cleared rather than set to a number.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20: Shut Up and Take My Gil
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 501.7143
$ws.Range("I33").Value = 349.92856
$ws.Range("J33").Value = 805.2857
$ws.Range("K33").Value = 349.92856
$ws.Range("L33").Value = 805.2857
$ws.Range("M33").Value = -120.92856
$ws.Range("N33").Value = -1263.2857
# Row 35: Conspicuous Conjuration
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 113: Amaro Kart
$ws.Range("H113").Value = 4692.4287
$ws.Range("I113").Value = 4193.1113
$ws.Range("K113").Value = 4193.1113
$ws.Range("M113").Value = -939.1112999999996
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 3419414.5
$ws.Range("I132").Value = 4223380
$ws.Range("K132").Value = 12670140
$ws.Range("M132").Value = -12667610
# Row 138: All-night Crafting
$ws.Range("H138").Value = 20573.639
$ws.Range("I138").Value = 2258.75
$ws.Range("J138").Value = 37667.535
$ws.Range("K138").Value = 6776.25
$ws.Range("L138").Value = 113002.605
$ws.Range("M138").Value = -1636.25
$ws.Range("N138").Value = -123282.605

$ws = $wb.Worksheets.Item("ARM")
# Row 63: Rivets Run through It
$ws.Range("H63").Value = 2690.4546
$ws.Range("I63").Value = 1920
$ws.Range("K63").Value = 1920
$ws.Range("M63").Value = -1234
# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 2690.4546
$ws.Range("I66").Value = 1920
$ws.Range("K66").Value = 9600
$ws.Range("M66").Value = -6168
# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1252.1072
$ws.Range("I110").Value = 917.5714
$ws.Range("K110").Value = 917.5714
$ws.Range("M110").Value = 1127.4286
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1722.3334
$ws.Range("I132").Value = 1454.7727
$ws.Range("J132").Value = 2142.7856
$ws.Range("K132").Value = 4364.3181
$ws.Range("L132").Value = 6428.3568
$ws.Range("M132").Value = -1834.3181
$ws.Range("N132").Value = -11488.3568

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 27703.8
$ws.Range("I20").Value = 41822.31
$ws.Range("J20").Value = 1483.7142
$ws.Range("K20").Value = 41822.31
$ws.Range("L20").Value = 1483.7142
$ws.Range("M20").Value = -41575.31
$ws.Range("N20").Value = -1977.7142
# Row 80: Unbreaker
$ws.Range("H80").Value = 602.3043
$ws.Range("J80").Value = 437.5
$ws.Range("L80").Value = 437.5
$ws.Range("N80").Value = -2433.5
# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 602.3043
$ws.Range("J83").Value = 437.5
$ws.Range("L83").Value = 2187.5
$ws.Range("N83").Value = -12171.5
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2158.1667
$ws.Range("I86").Value = 2158.1667
$ws.Range("K86").Value = 2158.1667
$ws.Range("M86").Value = -1035.1667
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2158.1667
$ws.Range("I89").Value = 2158.1667
$ws.Range("K89").Value = 10790.8335
$ws.Range("M89").Value = -5174.833500000001
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2336.5454
$ws.Range("I105").Value = 2078.2222
$ws.Range("J105").Value = 3499
$ws.Range("K105").Value = 2078.2222
$ws.Range("L105").Value = 3499
$ws.Range("M105").Value = -331.2222000000002
$ws.Range("N105").Value = -6993
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 2611.6428
$ws.Range("I134").Value = 2154.1936
$ws.Range("K134").Value = 6462.5808
$ws.Range("M134").Value = -3927.5808

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2780437.8
$ws.Range("I31").Value = 4546480
$ws.Range("J31").Value = 5228.9287
$ws.Range("K31").Value = 4546480
$ws.Range("L31").Value = 5228.9287
$ws.Range("M31").Value = -4546185
$ws.Range("N31").Value = -5818.9287
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2780437.8
$ws.Range("I34").Value = 4546480
$ws.Range("J34").Value = 5228.9287
$ws.Range("K34").Value = 4546480
$ws.Range("L34").Value = 5228.9287
$ws.Range("M34").Value = -4546278
$ws.Range("N34").Value = -5632.9287
# Row 135: The Wing's Wings
$ws.Range("H135").Value = 88609
$ws.Range("J135").Value = 88609
$ws.Range("L135").Value = 88609
$ws.Range("N135").Value = -98749

$ws = $wb.Worksheets.Item("CUL")
# Row 34: Fever Pitch
$ws.Range("H34").Value = 8157
$ws.Range("I34").Value = 849.75
$ws.Range("K34").Value = 2549.25
$ws.Range("M34").Value = -2465.25
# Row 50: Moving Up in the World
$ws.Range("H50").Value = 193.26315
$ws.Range("J50").Value = 198.16667
$ws.Range("L50").Value = 594.50001
$ws.Range("N50").Value = -1556.50001
# Row 53: Rolanberry Fields Forever
$ws.Range("H53").Value = 193.26315
$ws.Range("J53").Value = 198.16667
$ws.Range("L53").Value = 594.50001
$ws.Range("N53").Value = -1556.50001
# Row 107: Slippery Service
$ws.Range("H107").Value = 2039.238
$ws.Range("J107").Value = 858.1667
$ws.Range("L107").Value = 2574.5001
$ws.Range("N107").Value = -6414.5001

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 12317.333
$ws.Range("I70").Value = 13684.556
$ws.Range("K70").Value = 13684.556
$ws.Range("M70").Value = -13414.556
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 12317.333
$ws.Range("I73").Value = 13684.556
$ws.Range("K73").Value = 13684.556
$ws.Range("M73").Value = -12748.556
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 6104.8965
$ws.Range("I102").Value = 7580.4
$ws.Range("J102").Value = 2826
$ws.Range("K102").Value = 7580.4
$ws.Range("L102").Value = 2826
$ws.Range("M102").Value = -5958.4
$ws.Range("N102").Value = -6070
# Row 141: Mask Maker
$ws.Range("H141").Value = 73904.39999999999
$ws.Range("J141").Value = 73904.39999999999
$ws.Range("L141").Value = 73904.39999999999
$ws.Range("N141").Value = -84264.39999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 2724.875
$ws.Range("I7").Value = 2685.5715
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2685.5715
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2573.5715
$ws.Range("N7").Value = -3224
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 1006.9
$ws.Range("I22").Value = 434.75
$ws.Range("K22").Value = 434.75
$ws.Range("M22").Value = -139.75
# Row 27: Fire and Hide
$ws.Range("H27").Value = 1006.9
$ws.Range("I27").Value = 434.75
$ws.Range("K27").Value = 434.75
$ws.Range("M27").Value = -327.75
# Row 32: Men Who Scare Up Goats
$ws.Range("H32").Value = 47333
$ws.Range("I32").Value = 2000
$ws.Range("K32").Value = 2000
$ws.Range("M32").Value = -1683
# Row 40: Best Served Toad
$ws.Range("H40").Value = 3180.1538
$ws.Range("I40").Value = 2492.7778
$ws.Range("K40").Value = 2492.7778
$ws.Range("M40").Value = -2356.7778
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 4551.1816
$ws.Range("I46").Value = 719.8
$ws.Range("J46").Value = 7744
$ws.Range("K46").Value = 719.8
$ws.Range("L46").Value = 7744
$ws.Range("M46").Value = -531.8
$ws.Range("N46").Value = -8120
# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 1162.2222
$ws.Range("I61").Value = 855
$ws.Range("K61").Value = 855
$ws.Range("M61").Value = -653
# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 3097
$ws.Range("I93").Value = 4500
$ws.Range("K93").Value = 4500
$ws.Range("M93").Value = -3252
# Row 113: Peace in Rest
$ws.Range("H113").Value = 1162.2222
$ws.Range("I113").Value = 855
$ws.Range("K113").Value = 855
$ws.Range("M113").Value = 1315
# Row 122: Hell on Leather
$ws.Range("H122").Value = 4919.32
$ws.Range("I122").Value = 4036.3333
$ws.Range("K122").Value = 12108.9999
$ws.Range("M122").Value = -9658.999899999999
# Row 126: Battered Books
$ws.Range("H126").Value = 2724.875
$ws.Range("I126").Value = 2685.5715
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8056.7145
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5586.7145
$ws.Range("N126").Value = -13940
# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 3328.0625
$ws.Range("I132").Value = 3433.182
$ws.Range("J132").Value = 3096.8
$ws.Range("K132").Value = 10299.546
$ws.Range("L132").Value = 9290.400000000001
$ws.Range("M132").Value = -7769.545999999998
$ws.Range("N132").Value = -14350.4
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 3055.9092
$ws.Range("I136").Value = 3055.9092
$ws.Range("K136").Value = 9167.7276
$ws.Range("M136").Value = -6617.7276

$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 24855.445
$ws.Range("I81").Value = 34832.832
$ws.Range("J81").Value = 4900.6665
$ws.Range("K81").Value = 69665.664
$ws.Range("L81").Value = 9801.333000000001
$ws.Range("M81").Value = -68604.664
$ws.Range("N81").Value = -11923.333
# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 24855.445
$ws.Range("I84").Value = 34832.832
$ws.Range("J84").Value = 4900.6665
$ws.Range("K84").Value = 348328.32
$ws.Range("L84").Value = 49006.665
$ws.Range("M84").Value = -343024.32
$ws.Range("N84").Value = -59614.665
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 19984.305
$ws.Range("I136").Value = 24561.643
$ws.Range("J136").Value = 3963.625
$ws.Range("K136").Value = 73684.929
$ws.Range("L136").Value = 11890.875
$ws.Range("M136").Value = -71134.929
$ws.Range("N136").Value = -16990.875
